# Automatische test-sync: 2025-06-20 15:30:50
# Add the new incoming-mail row to the "Logs" sheet, extend the
# conditional-formatting ranges to cover it, and refresh the category
# breakdown on the "Dashboard" sheet to account for the new row.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 17 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(17, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item(17, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(17, 3).Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Cells.Item(17, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item(17, 5).Value = "Beste [naam],`nBedankt voor je interesse. Wij zijn geopend van maandag tot en met vrijdag van 9:00 tot 18:00 uur en op zaterdag van 10:00 tot 15:00 uur. Op zondag zijn wij gesloten. Mocht je nog verdere vragen hebben, laat het ons gerust weten.`nMet vriendelijke groet,`n[Naam]"
$logs.Cells.Item(17, 6).Value = "2025-06-20 15:30:13"
$logs.Cells.Item(17, 7).Value = "Ja"

# Reset the auto-grown row height back to the sheet default so the new
# row matches the plain (no explicit height) rows above it.
$logs.Rows.Item(17).AutoFit()

# --- Extend the conditional-formatting ranges to include row 17 ----------
$catRule = $logs.Range("D2:D16").FormatConditions
$catRule.Item(1).ModifyAppliesToRange($logs.Range("D2:D17"))

$answeredRule = $logs.Range("G2:G16").FormatConditions
$answeredRule.Item(1).ModifyAppliesToRange($logs.Range("G2:G17"))

# --- Dashboard sheet: refresh the category counts / ordering -------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(5, 2).Value = 2
$dash.Cells.Item(6, 1).Value = "Productinformatie"
$dash.Cells.Item(6, 2).Value = 1
$dash.Cells.Item(7, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(7, 2).Value = 1
$dash.Cells.Item(8, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(8, 2).Value = 1
